# The two species records stored in rows 17 and 18 were swapped (same
# column layout, the values for each record moved to the other row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values actually differ between row 17 and row 18 and need to
# be exchanged. (Columns such as C, I, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AT, AY are identical between the two rows, so they are left alone.)
$cols = @("A", "B", "D", "E", "F", "G", "H", "P", "Q", "R", "AW", "AX")

foreach ($col in $cols) {
    $addr17 = $col + "17"
    $addr18 = $col + "18"
    $v17 = $ws.Range($addr17).Value()
    $v18 = $ws.Range($addr18).Value()
    $ws.Range($addr17).Value = $v18
    $ws.Range($addr18).Value = $v17
}

# The blank "Ålder-Stadium" (K) cell also belongs to the record and moves
# from row 18 to row 17 along with the rest of the data.
$ws.Range("K18").Copy($ws.Range("K17"))
$ws.Range("K18").ClearContents()
